# Generate Report for Handback
# Update the handback-status workbook with refreshed timestamps / status
# for this handback run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G) for the
# 584e0deb... and acde9c89... rows moved forward.
$wsOverview.Range("G2").Value = "2016-08-25 22:16:43"
$wsOverview.Range("G4").Value = "2016-08-25 22:16:43"

# zh-cn sheet: status changed from "ht" to "mt" for the rows that share
# that status, and the handoff/handback datetimes advanced.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-25 22:16:37"
$wsZhCn.Range("H4").Value = "2016-08-25 22:16:37"
$wsZhCn.Range("K2").Value = "2016-08-25 22:16:55"
$wsZhCn.Range("K4").Value = "2016-08-25 22:16:55"

# de-de sheet: same status change plus the correspond handoff datetime
# (column H) and correspond handback datetime (column K) advanced.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-25 22:16:43"
$wsDeDe.Range("H4").Value = "2016-08-25 22:16:43"
$wsDeDe.Range("K2").Value = "2016-08-25 22:17:07"
$wsDeDe.Range("K4").Value = "2016-08-25 22:17:07"
